# Added toggle button functionality.
# This adds a new "New Text<value>" translation row to the "Translation" sheet,
# mirroring the pattern used by the TouchGFX text editor when a new text is created.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$newRow = 29

$ws.Cells.Item($newRow, 2).Value = "SingleUseId26"
$ws.Cells.Item($newRow, 3).Value = "Default"
$ws.Cells.Item($newRow, 4).Value = "Center"
$ws.Cells.Item($newRow, 5).Value = "LTR"
$ws.Cells.Item($newRow, 6).Value = "New Text<value>"
